$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139 (shifts existing rows 139-156 down to 140-157)
$ws.Rows.Item(139).Insert()

# Populate the new row 139 with the new weekly record
$ws.Range("A139").Value = 11
$ws.Range("B139").Value = "Vega Monumental Concepción"
$ws.Range("C139").Value = "Bíobío"
$ws.Range("D139").Value = "2021-10-05"
$ws.Range("E139").Value = 8
$ws.Range("F139").Value = "Fruta"
$ws.Range("G139").Value = 100101
$ws.Range("H139").Value = "Berries"
$ws.Range("I139").Value = 100112025
$ws.Range("J139").Value = "Frutilla"
$ws.Range("K139").Value = "Sin especificar"
$ws.Range("L139").Value = "Especial"
$ws.Range("M139").Value = 100
$ws.Range("N139").Value = 15000
$ws.Range("O139").Value = 16000
$ws.Range("P139").Value = 15500
$ws.Range("Q139").Value = "$/bandeja 7 kilos"
$ws.Range("R139").Value = "Provincia de Melipilla"
$ws.Range("S139").Value = 2214
$ws.Range("T139").Value = 7
